# Add a new "Docker" worksheet at the end of the workbook with a small
# command/usage reference table (matching the author's "changed file
# extension to .txt" commit, which in this workbook bundled a new Docker
# cheat-sheet tab).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end
# of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$dockerSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$dockerSheet.Name = "Docker"

# Command / usage rows.
$dockerSheet.Range("A1").Value = 'docker logs "container-name"'
$dockerSheet.Range("B1").Value = "logs all messages related to that container"

$dockerSheet.Range("A2").Value = 'docker kill "container-name"'
$dockerSheet.Range("B2").Value = "Stops the container - stopped container"

$dockerSheet.Range("A3").Value = 'docker rm "container-name"'
$dockerSheet.Range("B3").Value = "Removes / Deletes the container"

# Fit the two columns to their content, like the original author did.
$dockerSheet.Columns.Item(1).AutoFit() | Out-Null
$dockerSheet.Columns.Item(2).AutoFit() | Out-Null

# Make the new sheet the active / selected tab, with B3 selected.
$dockerSheet.Activate()
$dockerSheet.Range("B3").Select() | Out-Null
